$d = $word.ActiveDocument

# Locate the paragraph that contains the "Not using chunked encoding..." sentence
# (the last line of the bullet about providing object sizes / avoiding chunked encoding)
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Not using chunked encoding can also help, as we can atleast know the size once we receive header (a lot of servers didnt send a content length and relied on chunked encoding)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate target paragraph text"
}

# Determine which paragraph in the document this range falls in
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $searchRange.Start -and $p.Range.End -ge $searchRange.End) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not determine paragraph index for found text"
}

$targetPara = $d.Paragraphs.Item($targetIndex)
$insertRange = $targetPara.Range
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()

# The newly created paragraph inherits the bullet/list formatting of $targetPara.
# Fill in its text.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Automatically send responses for dependencies of an object in the response of the object itself. However, it can result in inefficient use of cache by browsers."
